# Mise à jour du cahier des charges
# Updates the functional specification sheet:
#  - Fcont1 ("Avoir un connecteur audio d'entrée") criterion/level no longer
#    requires a guitar input and now targets a 3.5mm jack instead of a 1/4" jack
#  - Fcont7 ("Respecter le budget estimé") materials cost is filled in (330)
#  - Selection moved to C22 to reflect where the user was working

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("C21").Value = "Pouvoir accueillir un micro"
$ws.Range("D21").Value = "Jack 3,5mm"
$ws.Range("D28").Value = "Matériel : 330"

$ws.Range("C22").Select()
